# Weekly price-sheet update: a new daily price record is inserted as row 20
# (Femacal de La Calera / Arándano (blue), Provincia de Limarí, 2021-09-23),
# pushing the previously-existing rows 20-98 down to rows 21-99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 20; Excel shifts rows 20..98 -> 21..99
# and carries the column D date formatting down with them.
$ws.Rows.Item(20).Insert()

# Populate the newly-inserted row 20 with the new record.
$ws.Range('A20').Value = 3
$ws.Range('B20').Value = 'Femacal de La Calera'
$ws.Range('C20').Value = 'Coquimbo'
$ws.Range('D20').Value = 44462
$ws.Range('E20').Value = 5
$ws.Range('F20').Value = 'Fruta'
$ws.Range('G20').Value = 100101
$ws.Range('H20').Value = 'Berries'
$ws.Range('I20').Value = 100101001
$ws.Range('J20').Value = 'Arándano (blue)'
$ws.Range('K20').Value = 'Sin especificar'
$ws.Range('L20').Value = 'Primera'
$ws.Range('M20').Value = 45
$ws.Range('N20').Value = 11000
$ws.Range('O20').Value = 11000
$ws.Range('P20').Value = 11000
$ws.Range('Q20').Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range('R20').Value = 'Provincia de Limarí'
$ws.Range('S20').Value = 7333
$ws.Range('T20').Value = 1.5
